$d = $word.ActiveDocument

# --- Paragraph 1: plain text paragraph ---
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()
$pIntro = $d.Paragraphs.Last
$pIntro.Range.ListFormat.RemoveNumbers()
$pIntro.Style = "Normal"
$pIntro.Range.Text = "How to create a tab for the tutorials:"

# --- Paragraph 2: hyperlink to w3schools tabs tutorial ---
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()
$pLink1 = $d.Paragraphs.Last
$pLink1.Range.ListFormat.RemoveNumbers()
$pLink1.Style = "Normal"
$pLink1.Range.Text = "How To Create Tabs (w3schools.com)"
$anchor1 = $d.Range($pLink1.Range.Start, $pLink1.Range.Start + 35)
$d.Hyperlinks.Add($anchor1, "https://www.w3schools.com/howto/howto_js_tabs.asp", $null, $null, "How To Create Tabs (w3schools.com)") | Out-Null

# --- Paragraph 3: hyperlink to React Tabs examples ---
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()
$pLink2 = $d.Paragraphs.Last
$pLink2.Range.ListFormat.RemoveNumbers()
$pLink2.Style = "Normal"
$pLink2.Range.Text = "React Tabs Examples | React School"
$anchor2 = $d.Range($pLink2.Range.Start, $pLink2.Range.Start + 35)
$d.Hyperlinks.Add($anchor2, "https://reactschool.us/tabs", $null, $null, "React Tabs Examples | React School") | Out-Null

Write-Output ("Done. ParaCount=" + $d.Paragraphs.Count)
